$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Chirayu"
$ws.Range("B1").Value = "+91 72489 46823"
$ws.Range("C1").Value = "qsda"
$ws.Range("D1").Value = "Dr. Sarah Davis: Pediatrician"
$ws.Range("E1").NumberFormat = "@"
$ws.Range("E1").Value = "2024-09-12"
$ws.Range("F1").Value = "14:12"
